$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 911.7564
$ws.Range("J17").Value = 851.0635
$ws.Range("L17").Value = 2553.1905
$ws.Range("N17").Value = -2889.1905

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 653.8276
$ws.Range("I28").Value = 646.2381
$ws.Range("J28").Value = 673.75
$ws.Range("K28").Value = 646.2381
$ws.Range("L28").Value = 673.75
$ws.Range("M28").Value = -161.2381
$ws.Range("N28").Value = -1643.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50924

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 47627784
$ws.Range("I132").Value = 58832384
$ws.Range("K132").Value = 176497152
$ws.Range("M132").Value = -176494622

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2568.851
$ws.Range("I137").Value = 1176.1154
$ws.Range("J137").Value = 4293.1904
$ws.Range("K137").Value = 3528.3462
$ws.Range("L137").Value = 12879.5712
$ws.Range("M137").Value = -978.3462
$ws.Range("N137").Value = -17979.5712

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5488.71
$ws.Range("I138").Value = 1074.5333
$ws.Range("J138").Value = 6267.682
$ws.Range("K138").Value = 3223.5999
$ws.Range("L138").Value = 18803.046
$ws.Range("M138").Value = 1916.4001
$ws.Range("N138").Value = -29083.046

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3421.2437
$ws.Range("I32").Value = 3071.8088
$ws.Range("J32").Value = 5797.4
$ws.Range("K32").Value = 3071.8088
$ws.Range("L32").Value = 5797.4
$ws.Range("M32").Value = -2784.8088
$ws.Range("N32").Value = -6371.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1262.5
$ws.Range("I45").Value = 1200
$ws.Range("K45").Value = 1200
$ws.Range("M45").Value = -823

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1109.9445
$ws.Range("I61").Value = 841.6445
$ws.Range("J61").Value = 2451.4443
$ws.Range("K61").Value = 841.6445
$ws.Range("L61").Value = 2451.4443
$ws.Range("M61").Value = -629.6445
$ws.Range("N61").Value = -2875.4443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3047.1162
$ws.Range("I74").Value = 2939.389
$ws.Range("J74").Value = 3601.1428
$ws.Range("K74").Value = 2939.389
$ws.Range("L74").Value = 3601.1428
$ws.Range("M74").Value = -2065.389
$ws.Range("N74").Value = -5349.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3047.1162
$ws.Range("I77").Value = 2939.389
$ws.Range("J77").Value = 3601.1428
$ws.Range("K77").Value = 14696.945
$ws.Range("L77").Value = 18005.714
$ws.Range("M77").Value = -10328.945
$ws.Range("N77").Value = -26741.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2089.5667
$ws.Range("I122").Value = 1434.75
$ws.Range("K122").Value = 4304.25
$ws.Range("M122").Value = -1854.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1992.3684
$ws.Range("I132").Value = 1173.2858
$ws.Range("K132").Value = 3519.8574
$ws.Range("M132").Value = -989.8574000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1109.9445
$ws.Range("I136").Value = 841.6445
$ws.Range("J136").Value = 2451.4443
$ws.Range("K136").Value = 2524.9335
$ws.Range("L136").Value = 7354.3329
$ws.Range("M136").Value = 25.06649999999991
$ws.Range("N136").Value = -12454.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 36437.5
$ws.Range("J137").Value = 39785.715
$ws.Range("L137").Value = 39785.715
$ws.Range("N137").Value = -49985.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 54500
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 54500
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 54500
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -55088

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 31410
$ws.Range("J60").Value = 31410
$ws.Range("L60").Value = 31410
$ws.Range("N60").Value = -32608

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1715.1279
$ws.Range("I134").Value = 1072.0968
$ws.Range("J134").Value = 3376.2917
$ws.Range("K134").Value = 3216.2904
$ws.Range("L134").Value = 10128.8751
$ws.Range("M134").Value = -681.2903999999999
$ws.Range("N134").Value = -15198.8751

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 44457
$ws.Range("J137").Value = 49321.25
$ws.Range("L137").Value = 49321.25
$ws.Range("N137").Value = -59521.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7465111.5
$ws.Range("I31").Value = 1314.234
$ws.Range("J31").Value = 25005036
$ws.Range("K31").Value = 1314.234
$ws.Range("L31").Value = 25005036
$ws.Range("M31").Value = -1019.234
$ws.Range("N31").Value = -25005626

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7465111.5
$ws.Range("I34").Value = 1314.234
$ws.Range("J34").Value = 25005036
$ws.Range("K34").Value = 1314.234
$ws.Range("L34").Value = 25005036
$ws.Range("M34").Value = -1112.234
$ws.Range("N34").Value = -25005440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1871.075
$ws.Range("I58").Value = 1807
$ws.Range("J58").Value = 2012.04
$ws.Range("K58").Value = 1807
$ws.Range("L58").Value = 2012.04
$ws.Range("M58").Value = -1604
$ws.Range("N58").Value = -2418.04

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2720.7856
$ws.Range("I122").Value = 1613.8572
$ws.Range("J122").Value = 3827.7144
$ws.Range("K122").Value = 4841.571599999999
$ws.Range("L122").Value = 11483.1432
$ws.Range("M122").Value = -2391.571599999999
$ws.Range("N122").Value = -16383.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3726.0833
$ws.Range("I134").Value = 4401.857
$ws.Range("J134").Value = 2780
$ws.Range("K134").Value = 13205.571
$ws.Range("L134").Value = 8340
$ws.Range("M134").Value = -10670.571
$ws.Range("N134").Value = -13410

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1871.075
$ws.Range("I136").Value = 1807
$ws.Range("J136").Value = 2012.04
$ws.Range("K136").Value = 5421
$ws.Range("L136").Value = 6036.12
$ws.Range("M136").Value = -2871
$ws.Range("N136").Value = -11136.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1038.4615
$ws.Range("I17").Value = 1038.4615
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 3115.3845
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -2946.3845
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 12851.294
$ws.Range("I34").Value = 18464
$ws.Range("J34").Value = 9789.817999999999
$ws.Range("K34").Value = 55392
$ws.Range("L34").Value = 29369.454
$ws.Range("M34").Value = -55308
$ws.Range("N34").Value = -29537.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 13246.909
$ws.Range("J39").Value = 13771.6
$ws.Range("L39").Value = 41314.8
$ws.Range("N39").Value = -41902.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4516.4287
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 4935.8335
$ws.Range("K55").Value = 6000
$ws.Range("L55").Value = 14807.5005
$ws.Range("M55").Value = -5823
$ws.Range("N55").Value = -15161.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 3505.182
$ws.Range("J136").Value = 4233.125
$ws.Range("L136").Value = 12699.375
$ws.Range("N136").Value = -22899.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2480.7144
$ws.Range("I138").Value = 1873
$ws.Range("K138").Value = 5619
$ws.Range("M138").Value = -479

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3173.238
$ws.Range("I122").Value = 1494.4286
$ws.Range("J122").Value = 6530.857
$ws.Range("K122").Value = 4483.2858
$ws.Range("L122").Value = 19592.571
$ws.Range("M122").Value = -2033.2858
$ws.Range("N122").Value = -24492.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 42794
$ws.Range("J137").Value = 42794
$ws.Range("L137").Value = 42794
$ws.Range("N137").Value = -52994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 88111.60000000001
$ws.Range("J139").Value = 88111.60000000001
$ws.Range("L139").Value = 88111.60000000001
$ws.Range("N139").Value = -98391.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6049.6665
$ws.Range("I40").Value = 3973.9375
$ws.Range("J40").Value = 9068.909
$ws.Range("K40").Value = 3973.9375
$ws.Range("L40").Value = 9068.909
$ws.Range("M40").Value = -3837.9375
$ws.Range("N40").Value = -9340.909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3084.7
$ws.Range("I132").Value = 1092.6482
$ws.Range("J132").Value = 7222.0386
$ws.Range("K132").Value = 3277.9446
$ws.Range("L132").Value = 21666.1158
$ws.Range("M132").Value = -747.9446000000003
$ws.Range("N132").Value = -26726.1158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2240.0613
$ws.Range("I136").Value = 1254.075
$ws.Range("J136").Value = 6622.222
$ws.Range("K136").Value = 3762.225
$ws.Range("L136").Value = 19866.666
$ws.Range("M136").Value = -1212.225
$ws.Range("N136").Value = -24966.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5651077
$ws.Range("I132").Value = 550.9268
$ws.Range("J132").Value = 18521718
$ws.Range("K132").Value = 1652.7804
$ws.Range("L132").Value = 55565154
$ws.Range("M132").Value = 877.2196000000001
$ws.Range("N132").Value = -55570214

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1703.419
$ws.Range("I136").Value = 610.52
$ws.Range("J136").Value = 3980.2917
$ws.Range("K136").Value = 1831.56
$ws.Range("L136").Value = 11940.8751
$ws.Range("M136").Value = 718.4400000000001
$ws.Range("N136").Value = -17040.8751
